$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 10485.25
$ws.Range("I40").Value = 21780.2
$ws.Range("J40").Value = 2417.4285
$ws.Range("K40").Value = 21780.2
$ws.Range("L40").Value = 2417.4285
$ws.Range("M40").Value = -21605.2
$ws.Range("N40").Value = -2767.4285
$ws.Range("H44").Value = 45830
$ws.Range("J44").Value = 45830
$ws.Range("L44").Value = 45830
$ws.Range("N44").Value = -46754
$ws.Range("H94").Value = 651231.25
$ws.Range("I94").Value = 780887.7
$ws.Range("J94").Value = 2949
$ws.Range("K94").Value = 780887.7
$ws.Range("L94").Value = 2949
$ws.Range("M94").Value = -780436.7
$ws.Range("N94").Value = -3851
$ws.Range("H96").Value = 71508240
$ws.Range("I96").Value = 5063.2856
$ws.Range("J96").Value = 143011410
$ws.Range("K96").Value = 15189.8568
$ws.Range("L96").Value = 429034230
$ws.Range("M96").Value = -13816.8568
$ws.Range("N96").Value = -429036976
$ws.Range("H100").Value = 1182.3334
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H129").Value = 1315.3429
$ws.Range("I129").Value = 1934.7142
$ws.Range("J129").Value = 1160.5
$ws.Range("K129").Value = 5804.142599999999
$ws.Range("L129").Value = 3481.5
$ws.Range("M129").Value = -804.1425999999992
$ws.Range("N129").Value = -13481.5
$ws.Range("H137").Value = 4574.316
$ws.Range("I137").Value = 1316.9333
$ws.Range("J137").Value = 6698.696
$ws.Range("K137").Value = 3950.7999
$ws.Range("L137").Value = 20096.088
$ws.Range("M137").Value = -1400.7999
$ws.Range("N137").Value = -25196.088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9300
$ws.Range("I32").Value = 8271.343000000001
$ws.Range("J32").Value = 23701.2
$ws.Range("K32").Value = 8271.343000000001
$ws.Range("L32").Value = 23701.2
$ws.Range("M32").Value = -7984.343000000001
$ws.Range("N32").Value = -24275.2
$ws.Range("H61").Value = 1399.0714
$ws.Range("I61").Value = 1136
$ws.Range("J61").Value = 2240.9
$ws.Range("K61").Value = 1136
$ws.Range("L61").Value = 2240.9
$ws.Range("M61").Value = -924
$ws.Range("N61").Value = -2664.9
$ws.Range("H74").Value = 1174.0779
$ws.Range("I74").Value = 1063.3235
$ws.Range("J74").Value = 2010.8889
$ws.Range("K74").Value = 1063.3235
$ws.Range("L74").Value = 2010.8889
$ws.Range("M74").Value = -189.3235
$ws.Range("N74").Value = -3758.8889
$ws.Range("H77").Value = 1174.0779
$ws.Range("I77").Value = 1063.3235
$ws.Range("J77").Value = 2010.8889
$ws.Range("K77").Value = 5316.6175
$ws.Range("L77").Value = 10054.4445
$ws.Range("M77").Value = -948.6175000000003
$ws.Range("N77").Value = -18790.4445
$ws.Range("H102").Value = 11413.608
$ws.Range("I102").Value = 1742.75
$ws.Range("J102").Value = 21963.637
$ws.Range("K102").Value = 1742.75
$ws.Range("L102").Value = 21963.637
$ws.Range("M102").Value = -120.75
$ws.Range("N102").Value = -25207.637
$ws.Range("H110").Value = 1454.5807
$ws.Range("I110").Value = 1510.8214
$ws.Range("J110").Value = 929.6667
$ws.Range("K110").Value = 1510.8214
$ws.Range("L110").Value = 929.6667
$ws.Range("M110").Value = 534.1786
$ws.Range("N110").Value = -5019.6667
$ws.Range("H114").Value = 42305.668
$ws.Range("J114").Value = 42305.668
$ws.Range("L114").Value = 42305.668
$ws.Range("N114").Value = -50983.668
$ws.Range("H122").Value = 2509.4614
$ws.Range("I122").Value = 2602.0908
$ws.Range("K122").Value = 7806.2724
$ws.Range("M122").Value = -5356.2724
$ws.Range("H136").Value = 1399.0714
$ws.Range("I136").Value = 1136
$ws.Range("J136").Value = 2240.9
$ws.Range("K136").Value = 3408
$ws.Range("L136").Value = 6722.700000000001
$ws.Range("M136").Value = -858
$ws.Range("N136").Value = -11822.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H86").Value = 1804.875
$ws.Range("I86").Value = 1931.3334
$ws.Range("J86").Value = 1642.2858
$ws.Range("K86").Value = 1931.3334
$ws.Range("L86").Value = 1642.2858
$ws.Range("M86").Value = -808.3334
$ws.Range("N86").Value = -3888.2858
$ws.Range("H89").Value = 1804.875
$ws.Range("I89").Value = 1931.3334
$ws.Range("J89").Value = 1642.2858
$ws.Range("K89").Value = 9656.666999999999
$ws.Range("L89").Value = 8211.429
$ws.Range("M89").Value = -4040.666999999999
$ws.Range("N89").Value = -19443.429
$ws.Range("H126").Value = 47850.332
$ws.Range("J126").Value = 47850.332
$ws.Range("L126").Value = 47850.332
$ws.Range("N126").Value = -57730.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4972.161
$ws.Range("I31").Value = 3950.1667
$ws.Range("J31").Value = 5047.8643
$ws.Range("K31").Value = 3950.1667
$ws.Range("L31").Value = 5047.8643
$ws.Range("M31").Value = -3655.1667
$ws.Range("N31").Value = -5637.8643
$ws.Range("H34").Value = 4972.161
$ws.Range("I34").Value = 3950.1667
$ws.Range("J34").Value = 5047.8643
$ws.Range("K34").Value = 3950.1667
$ws.Range("L34").Value = 5047.8643
$ws.Range("M34").Value = -3748.1667
$ws.Range("N34").Value = -5451.8643
$ws.Range("H122").Value = 120676.2
$ws.Range("I122").Value = 172023.14
$ws.Range("K122").Value = 516069.42
$ws.Range("M122").Value = -513619.42
$ws.Range("H132").Value = 28253.057
$ws.Range("I132").Value = 1163.4474
$ws.Range("J132").Value = 96880.07000000001
$ws.Range("K132").Value = 3490.3422
$ws.Range("L132").Value = 290640.21
$ws.Range("M132").Value = -960.3422
$ws.Range("N132").Value = -295700.21

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 383.9565
$ws.Range("J23").Value = 383.22223
$ws.Range("L23").Value = 1149.66669
$ws.Range("N23").Value = -1619.66669
$ws.Range("H131").Value = 4548
$ws.Range("I131").Value = 33699.668
$ws.Range("J131").Value = 1532.3103
$ws.Range("K131").Value = 101099.004
$ws.Range("L131").Value = 4596.9309
$ws.Range("M131").Value = -96059.00399999999
$ws.Range("N131").Value = -14676.9309

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 17836
$ws.Range("J15").Value = 17836
$ws.Range("L15").Value = 17836
$ws.Range("N15").Value = -18412
$ws.Range("H81").Value = 17836
$ws.Range("J81").Value = 17836
$ws.Range("L81").Value = 17836
$ws.Range("N81").Value = -19832
$ws.Range("H84").Value = 17836
$ws.Range("J84").Value = 17836
$ws.Range("L84").Value = 53508
$ws.Range("N84").Value = -63492
$ws.Range("H111").Value = 28468.25
$ws.Range("J111").Value = 28468.25
$ws.Range("L111").Value = 28468.25
$ws.Range("N111").Value = -34602.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7384
$ws.Range("I46").Value = 4248.3335
$ws.Range("J46").Value = 10071.714
$ws.Range("K46").Value = 4248.3335
$ws.Range("L46").Value = 10071.714
$ws.Range("M46").Value = -4060.3335
$ws.Range("N46").Value = -10447.714
$ws.Range("H55").Value = 465.625
$ws.Range("I55").Value = 332
$ws.Range("J55").Value = 637.4286
$ws.Range("K55").Value = 332
$ws.Range("L55").Value = 637.4286
$ws.Range("M55").Value = -159
$ws.Range("N55").Value = -983.4286
$ws.Range("H93").Value = 1656.2106
$ws.Range("I93").Value = 997.55554
$ws.Range("J93").Value = 2249
$ws.Range("K93").Value = 997.55554
$ws.Range("L93").Value = 2249
$ws.Range("M93").Value = 250.44446
$ws.Range("N93").Value = -4745
$ws.Range("H111").Value = 46249
$ws.Range("J111").Value = 46249
$ws.Range("L111").Value = 46249
$ws.Range("N111").Value = -54429
$ws.Range("H133").Value = 39997.332
$ws.Range("J133").Value = 39997.332
$ws.Range("L133").Value = 39997.332
$ws.Range("N133").Value = -45057.332
$ws.Range("H136").Value = 1107.017
$ws.Range("I136").Value = 879.2766
$ws.Range("J136").Value = 1999
$ws.Range("K136").Value = 2637.8298
$ws.Range("L136").Value = 5997
$ws.Range("M136").Value = -87.82979999999998
$ws.Range("N136").Value = -11097
$ws.Range("H139").Value = 50166.332
$ws.Range("J139").Value = 50166.332
$ws.Range("L139").Value = 50166.332
$ws.Range("N139").Value = -60446.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1450
$ws.Range("J96").Value = 1450
$ws.Range("L96").Value = 1450
$ws.Range("N96").Value = -4196
$ws.Range("H119").Value = 40674
$ws.Range("J119").Value = 40674
$ws.Range("L119").Value = 40674
$ws.Range("N119").Value = -50350
$ws.Range("H136").Value = 298360.9
$ws.Range("I136").Value = 408129.75
$ws.Range("J136").Value = 1985.1
$ws.Range("K136").Value = 1224389.25
$ws.Range("L136").Value = 5955.299999999999
$ws.Range("M136").Value = -1221839.25
$ws.Range("N136").Value = -11055.3
